$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.101.95"
$ws.Range("E2").Value = "  -0.84%  "
$ws.Range("D3").Value = "1.674.49"
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").Value = "211.83"
$ws.Range("E5").Value = "  -2.96%  "
$ws.Range("D6").Value = "0.5250"
$ws.Range("E6").Value = "  -5.13%  "
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("D8").Value = "0.2672"
$ws.Range("E8").Value = "  -1.06%  "
$ws.Range("D9").Value = "0.06293"
$ws.Range("E9").Value = "  -3.20%  "
$ws.Range("D10").Value = "21.26"
$ws.Range("E10").Value = "  -3.87%  "
$ws.Range("D11").Value = "0.07600"
$ws.Range("E11").Value = "  +0.56%  "
$ws.Range("D12").Value = "1.686.15"
$ws.Range("E12").Value = "  +0.45%  "
$ws.Range("E13").Value = "  -0.97%  "
$ws.Range("D14").Value = "0.5676"
$ws.Range("E14").Value = "  -2.32%  "
$ws.Range("D15").Value = "0.000008134"
$ws.Range("E15").Value = "  -3.85%  "
$ws.Range("D16").Value = "65.64"
$ws.Range("E16").Value = "  +0.71%  "
$ws.Range("D17").Value = "26.120.25"
$ws.Range("E17").Value = "  -1.01%  "
$ws.Range("E18").Value = "  -0.17%  "
$ws.Range("D19").Value = "4.830"
$ws.Range("E19").Value = "  -2.21%  "
$ws.Range("D20").Value = "10.59"
$ws.Range("E20").Value = "  -2.87%  "
$ws.Range("D21").Value = "188.47"
$ws.Range("E21").Value = "  -1.36%  "
$ws.Range("D22").Value = "6.183"
$ws.Range("E22").Value = "  -0.81%  "
$ws.Range("D23").Value = "1.007"
$ws.Range("E23").Value = "  -0.17%  "
$ws.Range("D24").Value = "148.58"
$ws.Range("E24").Value = "  +0.72%  "
$ws.Range("D25").Value = "0.1251"
$ws.Range("E25").Value = "  -5.82%  "
$ws.Range("D26").Value = "7.630"
$ws.Range("E26").Value = "  -3.35%  "
$ws.Range("D27").Value = "15.76"
$ws.Range("E27").Value = "  -0.26%  "
$ws.Range("D28").Value = "0.06357"
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("E29").Value = "  -2.81%  "
$ws.Range("D30").Value = "1.291"
$ws.Range("E30").Value = "  -2.42%  "
$ws.Range("D31").Value = "3.539"
$ws.Range("E31").Value = "  -1.51%  "
$ws.Range("D32").Value = "3.530"
$ws.Range("E32").Value = "  -1.52%  "
$ws.Range("E33").Value = "  -0.54%  "
$ws.Range("E34").Value = "  -2.90%  "
$ws.Range("D35").Value = "2.417"
$ws.Range("E35").Value = "  +0.68%  "
$ws.Range("D36").Value = "0.6026"
$ws.Range("E36").Value = "  -3.20%  "
$ws.Range("D37").Value = "2.711"
$ws.Range("E37").Value = "  -0.29%  "
$ws.Range("D38").Value = "6.135"
$ws.Range("E38").Value = "  -1.66%  "
$ws.Range("D39").Value = "0.01619"
$ws.Range("E39").Value = "  -0.49%  "
$ws.Range("D40").Value = "1.091.74"
$ws.Range("E40").Value = "  -1.88%  "
$ws.Range("D41").Value = "0.8709"
$ws.Range("E41").Value = "  -0.49%  "
$ws.Range("E42").Value = "  -1.00%  "
$ws.Range("D43").Value = "99.99"
$ws.Range("E43").Value = "  -0.84%  "
$ws.Range("D44").Value = "1.828.07"
$ws.Range("E44").Value = "  -0.25%  "
$ws.Range("E45").Value = "  -0.03%  "
$ws.Range("E46").Value = "  -0.64%  "
$ws.Range("D47").Value = "1.003"
$ws.Range("E47").Value = "  -0.10%  "
$ws.Range("D48").Value = "0.05243"
$ws.Range("E48").Value = "  -0.62%  "
$ws.Range("D49").Value = "7.949"
$ws.Range("E49").Value = "  -3.16%  "
$ws.Range("D50").Value = "0.4273"
$ws.Range("E50").Value = "  -0.48%  "
$ws.Range("D51").Value = "5.941"
$ws.Range("E51").Value = "  -2.23%  "
